$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K" header) values are being regenerated (Strike# -> K).
# Update the computed K values for rows 2-9.
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 8
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 7
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 6
$ws.Range("G9").Value = 1
